# Fixed Error on saving all records for each subcategory on the category csv file
# Fixed Category XLS listing first records as header
#
# Changes:
#  1) pescados-y-mariscos: swap "Elaborados" and "Sushi" rows
#  2) platos-preparados: swap "Verduras Congeladas" and "Apanados" rows
#  3) cuidado-personal: swap "Papel Higiénico" and "Cuidado De Pies" rows
#  4) limpieza-de-cocina: append a new "Toallas y Paños De Cocina" row

$wb = $excel.ActiveWorkbook

function Swap-Rows {
    param($SheetName, $Row1, $Row2)
    $ws = $wb.Worksheets.Item($SheetName)
    $cols = 4
    for ($c = 1; $c -le $cols; $c++) {
        $v1 = $ws.Cells.Item($Row1, $c).Value()
        $v2 = $ws.Cells.Item($Row2, $c).Value()
        $ws.Cells.Item($Row1, $c).Value = $v2
        $ws.Cells.Item($Row2, $c).Value = $v1
    }
}

# 1) pescados-y-mariscos: row 3 = Elaborados, row 4 = Sushi -> swap
Swap-Rows "pescados-y-mariscos" 3 4

# 2) platos-preparados: row 9 = Verduras Congeladas, row 10 = Apanados -> swap
Swap-Rows "platos-preparados" 9 10

# 3) cuidado-personal: row 8 = Papel Higiénico, row 9 = Cuidado De Pies -> swap
Swap-Rows "cuidado-personal" 8 9

# 4) limpieza-de-cocina: append new row 7 with the "Toallas y Paños De Cocina" subcategory
$wsLimpieza = $wb.Worksheets.Item("limpieza-de-cocina")
$wsLimpieza.Cells.Item(7, 1).Value = "Toallas y Paños De Cocina"
$wsLimpieza.Cells.Item(7, 2).Value = "https://www.tiendasjumbo.co/supermercado/limpieza-de-cocina/toallas-y-panos-de-cocina?initialMap=c,c&initialQuery=supermercado/limpieza-de-cocina&map=category-1,category-2,category-3&order=OrderByNameASC"
$wsLimpieza.Cells.Item(7, 3).Value = $wsLimpieza.Cells.Item(1, 3).Value()
$wsLimpieza.Cells.Item(7, 4).Value = "toallas-y-panos-de-cocina"
